$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the shared "sample number" text from E7760 to E7420.
# This value is shared by column G across rows 2-13 (same shared string).
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 7).Value = "E7420"
}

# 2. Convert the literal boolean FALSE in H2:H13 into a real =FALSE() formula.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 8).Formula = "=FALSE()"
}

# 3. Move the active selection from H2:H13 to G2:G13 (active cell G2).
$ws.Range("G2:G13").Select() | Out-Null
